# Auto-generated Excel COM-interop edit script
# Applies numeric updates described by the Omega_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 2479.3684
$ws.Cells.Item(19, 10).Value = 2726.125
$ws.Cells.Item(19, 12).Value = 2726.125
$ws.Cells.Item(19, 14).Value = -3076.125

$ws.Cells.Item(28, 8).Value = 244821.94
$ws.Cells.Item(28, 9).Value = 314509.94
$ws.Cells.Item(28, 11).Value = 314509.94
$ws.Cells.Item(28, 14).Value = -314024.94

$ws.Cells.Item(40, 8).Value = 14746.1875
$ws.Cells.Item(40, 9).Value = 2787.8
$ws.Cells.Item(40, 10).Value = 20181.818
$ws.Cells.Item(40, 11).Value = 2787.8
$ws.Cells.Item(40, 12).Value = 20181.818
$ws.Cells.Item(40, 13).Value = -2612.8
$ws.Cells.Item(40, 14).Value = -20531.818

$ws.Cells.Item(43, 8).Value = 18198.666
$ws.Cells.Item(43, 9).Value = 22299.5
$ws.Cells.Item(43, 10).Value = 9997
$ws.Cells.Item(43, 11).Value = 22299.5
$ws.Cells.Item(43, 12).Value = 9997
$ws.Cells.Item(43, 13).Value = -22230.5
$ws.Cells.Item(43, 14).Value = -10135

$ws.Cells.Item(51, 8).Value = 31616.666
$ws.Cells.Item(51, 9).Value = 9900
$ws.Cells.Item(51, 10).Value = 53333.332
$ws.Cells.Item(51, 11).Value = 9900
$ws.Cells.Item(51, 12).Value = 53333.332
$ws.Cells.Item(51, 13).Value = -9416
$ws.Cells.Item(51, 14).Value = -54301.332

$ws.Cells.Item(53, 8).Value = 532
$ws.Cells.Item(53, 9).Value = 498.5
$ws.Cells.Item(53, 10).Value = 599
$ws.Cells.Item(53, 11).Value = 498.5
$ws.Cells.Item(53, 12).Value = 599
$ws.Cells.Item(53, 13).Value = 138.5
$ws.Cells.Item(53, 14).Value = -1873

$ws.Cells.Item(64, 8).Value = 7790.706
$ws.Cells.Item(64, 9).Value = 6724.727
$ws.Cells.Item(64, 11).Value = 6724.727
$ws.Cells.Item(64, 13).Value = -6476.727

$ws.Cells.Item(67, 8).Value = 7790.706
$ws.Cells.Item(67, 9).Value = 6724.727
$ws.Cells.Item(67, 11).Value = 6724.727
$ws.Cells.Item(67, 13).Value = -5866.727

$ws.Cells.Item(88, 8).Value = 1616.3914
$ws.Cells.Item(88, 10).Value = 1557.5883
$ws.Cells.Item(88, 12).Value = 1557.5883
$ws.Cells.Item(88, 14).Value = -2369.5883

$ws.Cells.Item(91, 8).Value = 1616.3914
$ws.Cells.Item(91, 10).Value = 1557.5883
$ws.Cells.Item(91, 12).Value = 1557.5883
$ws.Cells.Item(91, 14).Value = -4365.588299999999

$ws.Cells.Item(94, 8).Value = 538.8333
$ws.Cells.Item(94, 9).Value = 556.6
$ws.Cells.Item(94, 11).Value = 556.6
$ws.Cells.Item(94, 13).Value = -105.6

$ws.Cells.Item(98, 8).Value = 1064.5883
$ws.Cells.Item(98, 9).Value = 1068.6875
$ws.Cells.Item(98, 11).Value = 1068.6875
$ws.Cells.Item(98, 13).Value = 429.3125

$ws.Cells.Item(101, 8).Value = 2029.4
$ws.Cells.Item(101, 9).Value = 749
$ws.Cells.Item(101, 11).Value = 2247
$ws.Cells.Item(101, 13).Value = -625

$ws.Cells.Item(112, 8).Value = 3955.5293
$ws.Cells.Item(112, 10).Value = 4282.933
$ws.Cells.Item(112, 12).Value = 12848.799
$ws.Cells.Item(112, 14).Value = -15064.799

$ws.Cells.Item(113, 8).Value = 4743.8
$ws.Cells.Item(113, 9).Value = 4616.6665
$ws.Cells.Item(113, 10).Value = 4798.2856
$ws.Cells.Item(113, 11).Value = 4616.6665
$ws.Cells.Item(113, 12).Value = 4798.2856
$ws.Cells.Item(113, 13).Value = -1362.6665
$ws.Cells.Item(113, 14).Value = -11306.2856

$ws.Cells.Item(115, 8).Value = 1587.9
$ws.Cells.Item(115, 9).Value = 1587.9
$ws.Cells.Item(115, 11).Value = 4763.700000000001
$ws.Cells.Item(115, 13).Value = -3196.700000000001

$ws.Cells.Item(116, 8).Value = 4525.4165
$ws.Cells.Item(116, 9).Value = 4544.4443
$ws.Cells.Item(116, 11).Value = 4544.4443
$ws.Cells.Item(116, 13).Value = -1102.4443

$ws.Cells.Item(122, 8).Value = 1064.5883
$ws.Cells.Item(122, 9).Value = 1068.6875
$ws.Cells.Item(122, 11).Value = 3206.0625
$ws.Cells.Item(122, 13).Value = -756.0625

$ws.Cells.Item(137, 8).Value = 2117.375
$ws.Cells.Item(137, 9).Value = 1542.4445
$ws.Cells.Item(137, 10).Value = 2856.5715
$ws.Cells.Item(137, 11).Value = 4627.333500000001
$ws.Cells.Item(137, 12).Value = 8569.7145
$ws.Cells.Item(137, 13).Value = -2077.333500000001
$ws.Cells.Item(137, 14).Value = -13669.7145

$ws.Cells.Item(138, 8).Value = 3325.9565
$ws.Cells.Item(138, 10).Value = 4668.8076
$ws.Cells.Item(138, 12).Value = 14006.4228
$ws.Cells.Item(138, 14).Value = -24286.4228

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 12000
$ws.Cells.Item(3, 10).Value = 12000
$ws.Cells.Item(3, 12).Value = 12000
$ws.Cells.Item(3, 14).Value = -12230

$ws.Cells.Item(32, 8).Value = 2287.1372
$ws.Cells.Item(32, 9).Value = 683.95746
$ws.Cells.Item(32, 10).Value = 21124.5
$ws.Cells.Item(32, 11).Value = 683.95746
$ws.Cells.Item(32, 12).Value = 21124.5
$ws.Cells.Item(32, 13).Value = -396.95746
$ws.Cells.Item(32, 14).Value = -21698.5

$ws.Cells.Item(45, 8).Value = 3344.5
$ws.Cells.Item(45, 9).Value = 2341.5
$ws.Cells.Item(45, 11).Value = 2341.5
$ws.Cells.Item(45, 13).Value = -1964.5

$ws.Cells.Item(97, 8).Value = 1719.75
$ws.Cells.Item(97, 9).Value = 602.3158
$ws.Cells.Item(97, 11).Value = 602.3158
$ws.Cells.Item(97, 13).Value = -106.3158

$ws.Cells.Item(122, 8).Value = 1657.1034
$ws.Cells.Item(122, 10).Value = 2107.8333
$ws.Cells.Item(122, 12).Value = 6323.499899999999
$ws.Cells.Item(122, 14).Value = -11223.4999

$ws.Cells.Item(132, 8).Value = 4392.6787
$ws.Cells.Item(132, 9).Value = 4493.3335
$ws.Cells.Item(132, 10).Value = 3788.75
$ws.Cells.Item(132, 11).Value = 13480.0005
$ws.Cells.Item(132, 12).Value = 11366.25
$ws.Cells.Item(132, 13).Value = -10950.0005
$ws.Cells.Item(132, 14).Value = -16426.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 3503.9614
$ws.Cells.Item(99, 9).Value = 3265.3914
$ws.Cells.Item(99, 11).Value = 3265.3914
$ws.Cells.Item(99, 13).Value = -1767.3914

$ws.Cells.Item(105, 8).Value = 1837.9
$ws.Cells.Item(105, 9).Value = 1797.625
$ws.Cells.Item(105, 11).Value = 1797.625
$ws.Cells.Item(105, 13).Value = -50.625

$ws.Cells.Item(107, 8).Value = 3369.5715
$ws.Cells.Item(107, 9).Value = 3117.4
$ws.Cells.Item(107, 11).Value = 3117.4
$ws.Cells.Item(107, 13).Value = -1197.4

$ws.Cells.Item(134, 8).Value = 2770.647
$ws.Cells.Item(134, 9).Value = 2740.2666
$ws.Cells.Item(134, 11).Value = 8220.799800000001
$ws.Cells.Item(134, 13).Value = -5685.799800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 490.15384
$ws.Cells.Item(22, 9).Value = 506.08334
$ws.Cells.Item(22, 10).Value = 299
$ws.Cells.Item(22, 11).Value = 506.08334
$ws.Cells.Item(22, 12).Value = 299
$ws.Cells.Item(22, 13).Value = -156.08334
$ws.Cells.Item(22, 14).Value = -999

$ws.Cells.Item(31, 8).Value = 12914
$ws.Cells.Item(31, 9).Value = 16296.777
$ws.Cells.Item(31, 11).Value = 16296.777
$ws.Cells.Item(31, 13).Value = -16001.777

$ws.Cells.Item(34, 8).Value = 12914
$ws.Cells.Item(34, 9).Value = 16296.777
$ws.Cells.Item(34, 11).Value = 16296.777
$ws.Cells.Item(34, 13).Value = -16094.777

$ws.Cells.Item(55, 8).Value = 5857.6665
$ws.Cells.Item(55, 9).Value = 3786.5
$ws.Cells.Item(55, 10).Value = 10000
$ws.Cells.Item(55, 11).Value = 3786.5
$ws.Cells.Item(55, 12).Value = 10000
$ws.Cells.Item(55, 13).Value = -3471.5
$ws.Cells.Item(55, 14).Value = -10630

$ws.Cells.Item(86, 8).Value = 19614544
$ws.Cells.Item(86, 10).Value = 7782.125
$ws.Cells.Item(86, 12).Value = 7782.125
$ws.Cells.Item(86, 14).Value = -10028.125

$ws.Cells.Item(89, 8).Value = 19614544
$ws.Cells.Item(89, 10).Value = 7782.125
$ws.Cells.Item(89, 12).Value = 38910.625
$ws.Cells.Item(89, 14).Value = -50142.625

$ws.Cells.Item(99, 8).Value = 3476840
$ws.Cells.Item(99, 9).Value = 4634303.5
$ws.Cells.Item(99, 11).Value = 4634303.5
$ws.Cells.Item(99, 13).Value = -4632805.5

$ws.Cells.Item(107, 8).Value = 1065.1482
$ws.Cells.Item(107, 9).Value = 1054.4
$ws.Cells.Item(107, 11).Value = 1054.4
$ws.Cells.Item(107, 13).Value = 865.5999999999999

$ws.Cells.Item(126, 8).Value = 3476840
$ws.Cells.Item(126, 9).Value = 4634303.5
$ws.Cells.Item(126, 11).Value = 13902910.5
$ws.Cells.Item(126, 13).Value = -13900440.5

$ws.Cells.Item(132, 8).Value = 5096.16
$ws.Cells.Item(132, 9).Value = 4806.476
$ws.Cells.Item(132, 11).Value = 14419.428
$ws.Cells.Item(132, 13).Value = -11889.428

$ws.Cells.Item(134, 8).Value = 6017.357
$ws.Cells.Item(134, 10).Value = 4999
$ws.Cells.Item(134, 12).Value = 14997
$ws.Cells.Item(134, 14).Value = -20067

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1983
$ws.Cells.Item(5, 10).Value = 888
$ws.Cells.Item(5, 12).Value = 2664
$ws.Cells.Item(5, 14).Value = -2888

$ws.Cells.Item(33, 8).Value = 84.5
$ws.Cells.Item(33, 9).Value = 91.40000000000001
$ws.Cells.Item(33, 11).Value = 548.4000000000001
$ws.Cells.Item(33, 13).Value = -265.4000000000001

$ws.Cells.Item(38, 8).Value = 163.46153
$ws.Cells.Item(38, 10).Value = 175.45454
$ws.Cells.Item(38, 12).Value = 526.3636200000001
$ws.Cells.Item(38, 14).Value = -1220.36362

$ws.Cells.Item(39, 8).Value = 4481.8667
$ws.Cells.Item(39, 9).Value = 2611.6667
$ws.Cells.Item(39, 10).Value = 4949.4165
$ws.Cells.Item(39, 11).Value = 7835.000100000001
$ws.Cells.Item(39, 12).Value = 14848.2495
$ws.Cells.Item(39, 13).Value = -7541.000100000001
$ws.Cells.Item(39, 14).Value = -15436.2495

$ws.Cells.Item(68, 8).Value = 493.85715
$ws.Cells.Item(68, 9).Value = 373.4
$ws.Cells.Item(68, 10).Value = 795
$ws.Cells.Item(68, 11).Value = 1120.2
$ws.Cells.Item(68, 12).Value = 2385
$ws.Cells.Item(68, 13).Value = -309.1999999999998
$ws.Cells.Item(68, 14).Value = -4007

$ws.Cells.Item(71, 8).Value = 493.85715
$ws.Cells.Item(71, 9).Value = 373.4
$ws.Cells.Item(71, 10).Value = 795
$ws.Cells.Item(71, 11).Value = 3360.6
$ws.Cells.Item(71, 12).Value = 7155
$ws.Cells.Item(71, 13).Value = 695.4000000000001
$ws.Cells.Item(71, 14).Value = -15267

$ws.Cells.Item(118, 8).Value = 7998
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 12166.667
$ws.Cells.Item(122, 9).Value = 500
$ws.Cells.Item(122, 10).Value = 18000
$ws.Cells.Item(122, 11).Value = 4500
$ws.Cells.Item(122, 12).Value = 162000
$ws.Cells.Item(122, 13).Value = -2050
$ws.Cells.Item(122, 14).Value = -166900

$ws.Cells.Item(128, 8).Value = 165499.5
$ws.Cells.Item(128, 9).Value = 165499.5
$ws.Cells.Item(128, 11).Value = 496498.5
$ws.Cells.Item(128, 13).Value = -491518.5

$ws.Cells.Item(131, 8).Value = 2215.513
$ws.Cells.Item(131, 10).Value = 2399.375
$ws.Cells.Item(131, 12).Value = 7198.125
$ws.Cells.Item(131, 14).Value = -17278.125

$ws.Cells.Item(135, 8).Value = 1983
$ws.Cells.Item(135, 10).Value = 888
$ws.Cells.Item(135, 12).Value = 7992
$ws.Cells.Item(135, 14).Value = -13062

$ws.Cells.Item(137, 8).Value = 4773.25
$ws.Cells.Item(137, 9).Value = 2170.6
$ws.Cells.Item(137, 11).Value = 6511.799999999999
$ws.Cells.Item(137, 13).Value = -1411.799999999999

$ws.Cells.Item(141, 8).Value = 5639.1113
$ws.Cells.Item(141, 9).Value = 4393.4287
$ws.Cells.Item(141, 11).Value = 13180.2861
$ws.Cells.Item(141, 13).Value = -8000.286100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6000
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 14).ClearContents()

$ws.Cells.Item(73, 8).Value = 6000
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 14).ClearContents()

$ws.Cells.Item(107, 8).Value = 395.6
$ws.Cells.Item(107, 9).Value = 355.16666
$ws.Cells.Item(107, 11).Value = 355.16666
$ws.Cells.Item(107, 13).Value = 1564.83334

$ws.Cells.Item(132, 8).Value = 3681.8865
$ws.Cells.Item(132, 9).Value = 3418.1177
$ws.Cells.Item(132, 10).Value = 4578.7
$ws.Cells.Item(132, 11).Value = 10254.3531
$ws.Cells.Item(132, 12).Value = 13736.1
$ws.Cells.Item(132, 13).Value = -7724.3531
$ws.Cells.Item(132, 14).Value = -18796.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 14874.571
$ws.Cells.Item(7, 9).Value = 10628.75
$ws.Cells.Item(7, 11).Value = 10628.75
$ws.Cells.Item(7, 13).Value = -10516.75

$ws.Cells.Item(22, 8).Value = 2244.1
$ws.Cells.Item(22, 10).Value = 2341.4375
$ws.Cells.Item(22, 12).Value = 2341.4375
$ws.Cells.Item(22, 14).Value = -2931.4375

$ws.Cells.Item(27, 8).Value = 2244.1
$ws.Cells.Item(27, 10).Value = 2341.4375
$ws.Cells.Item(27, 12).Value = 2341.4375
$ws.Cells.Item(27, 14).Value = -2555.4375

$ws.Cells.Item(36, 8).Value = 69715
$ws.Cells.Item(36, 10).Value = 69715
$ws.Cells.Item(36, 12).Value = 69715
$ws.Cells.Item(36, 14).Value = -70839

$ws.Cells.Item(40, 8).Value = 11252.793
$ws.Cells.Item(40, 9).Value = 8946
$ws.Cells.Item(40, 10).Value = 17308.125
$ws.Cells.Item(40, 11).Value = 8946
$ws.Cells.Item(40, 12).Value = 17308.125
$ws.Cells.Item(40, 13).Value = -8810
$ws.Cells.Item(40, 14).Value = -17580.125

$ws.Cells.Item(46, 8).Value = 5055.2383
$ws.Cells.Item(46, 9).Value = 6582
$ws.Cells.Item(46, 11).Value = 6582
$ws.Cells.Item(46, 13).Value = -6394

$ws.Cells.Item(68, 8).Value = 9002
$ws.Cells.Item(68, 9).Value = 8000
$ws.Cells.Item(68, 11).Value = 8000
$ws.Cells.Item(68, 13).Value = -7251

$ws.Cells.Item(71, 8).Value = 9002
$ws.Cells.Item(71, 9).Value = 8000
$ws.Cells.Item(71, 11).Value = 40000
$ws.Cells.Item(71, 13).Value = -36256

$ws.Cells.Item(122, 8).Value = 5547.4443
$ws.Cells.Item(122, 9).Value = 5547.4443
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 16642.3329
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -14192.3329
$ws.Cells.Item(122, 14).ClearContents()

$ws.Cells.Item(126, 8).Value = 14874.571
$ws.Cells.Item(126, 9).Value = 10628.75
$ws.Cells.Item(126, 11).Value = 31886.25
$ws.Cells.Item(126, 13).Value = -29416.25

$ws.Cells.Item(132, 8).Value = 3567.4243
$ws.Cells.Item(132, 9).Value = 3513.4285
$ws.Cells.Item(132, 11).Value = 10540.2855
$ws.Cells.Item(132, 13).Value = -8010.2855

$ws.Cells.Item(133, 8).Value = 67499.5
$ws.Cells.Item(133, 10).Value = 67499.5
$ws.Cells.Item(133, 12).Value = 67499.5
$ws.Cells.Item(133, 14).Value = -72559.5

$ws.Cells.Item(134, 8).Value = 72999.5
$ws.Cells.Item(134, 10).Value = 72999.5
$ws.Cells.Item(134, 12).Value = 72999.5
$ws.Cells.Item(134, 14).Value = -83139.5

$ws.Cells.Item(136, 8).Value = 2199.8
$ws.Cells.Item(136, 9).Value = 2199.8
$ws.Cells.Item(136, 11).Value = 6599.400000000001
$ws.Cells.Item(136, 13).Value = -4049.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1984
$ws.Cells.Item(107, 9).Value = 1840.1666
$ws.Cells.Item(107, 10).Value = 2156.6
$ws.Cells.Item(107, 11).Value = 5520.4998
$ws.Cells.Item(107, 12).Value = 6469.799999999999
$ws.Cells.Item(107, 13).Value = -3600.4998
$ws.Cells.Item(107, 14).Value = -10309.8

$ws.Cells.Item(122, 8).Value = 5079.8945
$ws.Cells.Item(122, 9).Value = 4501.0586
$ws.Cells.Item(122, 10).Value = 10000
$ws.Cells.Item(122, 11).Value = 13503.1758
$ws.Cells.Item(122, 12).Value = 30000
$ws.Cells.Item(122, 13).Value = -11053.1758
$ws.Cells.Item(122, 14).Value = -34900

$ws.Cells.Item(132, 8).Value = 4073.6365
$ws.Cells.Item(132, 9).Value = 4077.889
$ws.Cells.Item(132, 10).Value = 4054.5
$ws.Cells.Item(132, 11).Value = 12233.667
$ws.Cells.Item(132, 12).Value = 12163.5
$ws.Cells.Item(132, 13).Value = -9703.667000000001
$ws.Cells.Item(132, 14).Value = -17223.5

$ws.Cells.Item(135, 8).Value = 84499.125
$ws.Cells.Item(135, 10).Value = 84499.125
$ws.Cells.Item(135, 12).Value = 84499.125
$ws.Cells.Item(135, 14).Value = -94639.125
